# Insert a new data row at row 32 (pushing existing rows 32-103 down to 33-104)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 32.
$ws.Rows("32").Insert()

# Populate the newly inserted row 32 with the new record.
$ws.Cells.Item(32, 1).Value2 = 5
$ws.Cells.Item(32, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(32, 3).Value2 = "Maule"
$ws.Cells.Item(32, 4).Value2 = 44608
$ws.Cells.Item(32, 5).Value2 = 7
$ws.Cells.Item(32, 6).Value2 = 100112001
$ws.Cells.Item(32, 7).Value2 = "Berenjena"
$ws.Cells.Item(32, 8).Value2 = "Sin especificar"
$ws.Cells.Item(32, 9).Value2 = "Primera"
$ws.Cells.Item(32, 10).Value2 = 200
$ws.Cells.Item(32, 11).Value2 = 7000
$ws.Cells.Item(32, 12).Value2 = 7000
$ws.Cells.Item(32, 13).Value2 = 7000
$ws.Cells.Item(32, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(32, 15).Value2 = "Región del Maule"
$ws.Cells.Item(32, 16).Value2 = 117
$ws.Cells.Item(32, 17).Value2 = 60
$ws.Cells.Item(32, 18).Value2 = "Hortaliza"

# Ensure the date cell keeps the same date number format as the rest of column D.
$ws.Cells.Item(32, 4).NumberFormat = $ws.Cells.Item(33, 4).NumberFormat
